$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the report title (A2) and source citation (A52): January 2017 -> February 2017
$ws.Range("A2").Value = "Short-Term Energy Outlook, February 2017"
$ws.Range("A52").Value = "Source: Short-Term Energy Outlook, February 2017."

# Updated Non-OPEC supply (C), World demand (D) and WTI price (E) underlying data.
# G/H/I (change-from-previous-year) columns are formulas and recompute automatically,
# and the chart series caches/external-link caches follow from those same cells.
$ws.Range("C28").Value = 53.536449724999997
$ws.Range("D28").Value = 91.167427353999997
$ws.Range("C29").Value = 54.174579627999996
$ws.Range("D29").Value = 91.851533747000005
$ws.Range("C30").Value = 55.061020433000003
$ws.Range("D30").Value = 93.036322092000006
$ws.Range("C31").Value = 55.894061153999999
$ws.Range("D31").Value = 93.217316842000002
$ws.Range("C32").Value = 55.832797691000003
$ws.Range("D32").Value = 92.872456141000001
$ws.Range("C33").Value = 56.802195955000002
$ws.Range("D33").Value = 92.945229963000003
$ws.Range("C34").Value = 57.423574254000002
$ws.Range("D34").Value = 93.985852453000007
$ws.Range("C35").Value = 58.669676477000003
$ws.Range("D35").Value = 94.590933238000005
$ws.Range("C36").Value = 58.453621707000003
$ws.Range("D36").Value = 94.057179778999995
$ws.Range("C37").Value = 58.498679959
$ws.Range("D37").Value = 94.595104372999998
$ws.Range("C38").Value = 58.994774317000001
$ws.Range("D38").Value = 96.034932707999999
$ws.Range("C39").Value = 59.225250119999998
$ws.Range("D39").Value = 95.526807508000005
$ws.Range("C40").Value = 58.508242992
$ws.Range("D40").Value = 95.360327831999996
$ws.Range("C41").Value = 57.708725993999998
$ws.Range("D41").Value = 96.087232920999995
$ws.Range("C42").Value = 57.978004562999999
$ws.Range("D42").Value = 97.460478365
$ws.Range("C43").Value = 58.594457708
$ws.Range("D43").Value = 96.940833552000001
$ws.Range("E43").Value = 49.179625000000001
$ws.Range("C44").Value = 57.706747063000002
$ws.Range("D44").Value = 96.983636528999995
$ws.Range("E44").Value = 52.83
$ws.Range("C45").Value = 58.329605641999997
$ws.Range("D45").Value = 97.723528970999993
$ws.Range("E45").Value = 53
$ws.Range("C46").Value = 58.734544735
$ws.Range("D46").Value = 98.995585095999999
$ws.Range("E46").Value = 54
$ws.Range("C47").Value = 59.127679077000003
$ws.Range("D47").Value = 98.621587374000001
$ws.Range("E47").Value = 54
$ws.Range("C48").Value = 58.730159952000001
$ws.Range("D48").Value = 98.578322159999999
$ws.Range("E48").Value = 54
$ws.Range("C49").Value = 59.537479816999998
$ws.Range("D49").Value = 99.180575309999995
$ws.Range("E49").Value = 55.671875
$ws.Range("C50").Value = 59.768345443999998
$ws.Range("D50").Value = 100.31969019
$ws.Range("E50").Value = 56.666666667000001
$ws.Range("C51").Value = 60.137244148000001
$ws.Range("D51").Value = 100.07765381
$ws.Range("E51").Value = 58.3125
